$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellFromBase64 {
    param($ws, $row, $b64)
    $bytes = [System.Convert]::FromBase64String($b64)
    $text = [System.Text.Encoding]::UTF8.GetString($bytes)
    $ws.Cells.Item($row, 5).Value = $text
}

Set-CellFromBase64 $ws 2 "W0ZlaSVaaG91JU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBUaW5nJVl1JU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBSb25naHVpJUR1JU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBHdW9odWklRmFuJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBZaW5nJUxpdSVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgWmhpYm8lTGl1JU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBKaWUlWGlhbmclTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIFllbWluZyVXYW5nJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBCaW4lU29uZyVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgWGlhb3lpbmclR3UlTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIEx1bHUlR3VhbiVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgWXVhbiVXZWklTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIEh1aSVMaSVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgWHVkb25nJVd1JU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBKaXV5YW5nJVh1JU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBTaGVuZ2ppbiVUdSVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgWWklWmhhbmclTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIEh1YSVDaGVuJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBCaW4lQ2FvJU5VTEwlMF0="
Set-CellFromBase64 $ws 3 "W0NhcmJcdTAwZjMtQ2lzbmVybyVZYWNxdWVsaW4lY29yZUdpdmVzTm9FbWFpbCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgRmVyblx1MDBlMW5kZXotR29uelx1MDBlMWxleiVQYXVsYSVjb3JlR2l2ZXNOb0VtYWlsJTAsICAgICAgICAgICAgICAgICAgICAgICAgICBIaWVycmV6dWVsby1Sb2phcyVOYWlmaSVjb3JlR2l2ZXNOb0VtYWlsJTAsICAgICAgICAgICAgICAgICAgICAgICAgICBTdWJlcnQtU2FsYXMlTGl6YW5kcmElY29yZUdpdmVzTm9FbWFpbCUwXQ=="
Set-CellFromBase64 $ws 4 "W0t5dW5nIFNvbyVIb25nJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBLd2FuIEhvJUxlZSVOVUxMJTIsICAgICAgICAgICAgICAgICAgICAgICAgICAgS3dhbiBIbyVMZWUlTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIEppbiBIb25nJUNodW5nJU5VTEwlMiwgICAgICAgICAgICAgICAgICAgICAgICAgICBKaW4gSG9uZyVDaHVuZyVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgS3llb25nLUNoZW9sJVNoaW4lTlVMTCUyLCAgICAgICAgICAgICAgICAgICAgICAgICAgIEt5ZW9uZy1DaGVvbCVTaGluJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBFdW4gWW91bmclQ2hvaSVOVUxMJTIsICAgICAgICAgICAgICAgICAgICAgICAgICAgRXVuIFlvdW5nJUNob2klTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIEh5dW4gSnVuZyVKaW4lTlVMTCUyLCAgICAgICAgICAgICAgICAgICAgICAgICAgIEh5dW4gSnVuZyVKaW4lTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIEpvbmcgR2VvbCVKYW5nJU5VTEwlMiwgICAgICAgICAgICAgICAgICAgICAgICAgICBKb25nIEdlb2wlSmFuZyVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgV29uaHdhJUxlZSVOVUxMJTIsICAgICAgICAgICAgICAgICAgICAgICAgICAgV29uaHdhJUxlZSVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgSnVuZSBIb25nJUFobiVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgSnVuZSBIb25nJUFobiVOVUxMJTBd"
Set-CellFromBase64 $ws 5 "W01hdHQlQXJlbnR6JXhyZWYgbm8gZW1haWwlMCwgICAgICAgRXJpYyVZaW0leHJlZiBubyBlbWFpbCUwLCAgICAgICBMaW5keSVLbGFmZiV4cmVmIG5vIGVtYWlsJTAsICAgICAgIFNoYXJ1a2glTG9raGFuZHdhbGEleHJlZiBubyBlbWFpbCUwLCAgICAgICBGcmFuY2lzIFguJVJpZWRvJXhyZWYgbm8gZW1haWwlMCwgICAgICAgTWFyaWElQ2hvbmcleHJlZiBubyBlbWFpbCUwLCAgICAgICBNZWxpc3NhJUxlZSV4cmVmIG5vIGVtYWlsJTBd"
Set-CellFromBase64 $ws 6 "W0FnZ2Fyd2FsJUdhdXJhdiVjb3JlR2l2ZXNOb0VtYWlsJTAsICAgICAgICAgICAgICAgICAgICAgICAgICBBZ2dhcndhbCVTYXVyYWJoJWNvcmVHaXZlc05vRW1haWwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgIEdhcmNpYS1UZWxsZXMlTmVsc29uJWNvcmVHaXZlc05vRW1haWwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgIEhlbnJ5JUJyYW5kb24gTWljaGFlbCVjb3JlR2l2ZXNOb0VtYWlsJTAsICAgICAgICAgICAgICAgICAgICAgICAgICBMYXZpZSVDYXJsJWNvcmVHaXZlc05vRW1haWwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgIExpcHBpJUdpdXNlcHBlJWNvcmVHaXZlc05vRW1haWwlMF0="
Set-CellFromBase64 $ws 7 "W1l1YW4lWXUlTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIERhbiVYdSVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgU2hvdXpoaSVGdSVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgSnVuJVpoYW5nJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBYaWFvYm8lWWFuZyVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgTGlhbmclWHUlTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIEppcWlhbiVYdSVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgWW9uZ3JhbiVXdSVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgQ2hhb2xpbiVIdWFuZyVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgWWFxaSVPdXlhbmclTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIEx1eXUlWWFuZyVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgTWluZ2hhbyVGYW5nJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBIb25nd2VuJVhpYW8lTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIEppbmclTWElTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIFdlaSVaaHUlTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIFNvbmclSHUlTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIFF1YW4lSHUlTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIERhb3lpbiVEaW5nJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBNaW5nJUh1JU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBHdW9jaGFvJVpodSVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgV2VpamlhbmclWHUlTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIEp1biVHdW8lTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIEppbmdsb25nJVh1JU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBIYWl0YW8lWXVhbiVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgQmluJVpoYW5nJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBaaHVpJVl1JXl1emh1aUB3aHUuZWR1LmNuJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgRGVjaGFuZyVDaGVuJWljdWRlY2hhbmdjaGVuQDE2My5jb20lMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBTaGl5aW5nJVl1YW4leXVhbl9zaGl5aW5nQDE2My5jb20lMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBZb3UlU2hhbmcleW91X3NoYW5naHVzdEAxNjMuY29tJTBd"
Set-CellFromBase64 $ws 8 "W0ZhbiVZYW5nJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBTaGFvYm8lU2hpJU5VTEwlMiwgICAgICAgICAgICAgICAgICAgICAgICAgICBTaGFvYm8lU2hpJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBKaWxpbmclWmh1JU5VTEwlMSwgICAgICAgICAgICAgICAgICAgICAgICAgICBKaW56aGklU2hpJU5VTEwlMSwgICAgICAgICAgICAgICAgICAgICAgICAgICBLYWklRGFpJU5VTEwlMSwgICAgICAgICAgICAgICAgICAgICAgICAgICBYaWFvYmVpJUNoZW4lMTIxMjc4Mzk1QHFxLmNvbSUxXQ=="
Set-CellFromBase64 $ws 9 "W0ppYS1GdSVXZWklTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIEZhbmctWWFuZyVIdWFuZyVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgVGlhbi1ZdWFuJVhpb25nJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBUaWFuLVl1YW4lWGlvbmclTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIFFpJUxpdSVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgSG9uZyVDaGVuJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBIdWklV2FuZyVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgSGUlSHVhbmclTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIFlpLUNodW4lTHVvJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBYdWFuJVpob3UlTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIFpoaS1ZdWUlTGl1JU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBZb25nJVBlbmclTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIFl1YW4tTmluZyVYdSVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgQm8lV2FuZyVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgWWluZy1ZaW5nJVlhbmclTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIFpvbmctQW4lTGlhbmclTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIFh1ZS1aaG9uZyVMZWklTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIFlhbmclR2UlTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIE1pbmclWWFuZyVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgTGluZyVaaGFuZyVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgTWluZy1RdWFuJVplbmclTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIEhlJVl1JU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBLYWklTGl1JU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBZdS1IZW5nJUppYSVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgQmVybmFyZCBEJVByZW5kZXJnYXN0JU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBXZWktTWluJUxpJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBXZWktTWluJUxpJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBNYW8lQ2hlbiVOVUxMJTBd"
Set-CellFromBase64 $ws 10 "W1hpYW9jaGVuJUxpJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBTaHV5dW4lWHUlTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIE11cWluZyVZdSVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgS2UlV2FuZyVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgWXUlVGFvJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBZaW5nJVpob3UlTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIEppbmclU2hpJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBNaW4lWmhvdSVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgQm8lV3UlTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIFpoZW55dSVZYW5nJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBDb25nJVpoYW5nJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBKdW5xaW5nJVl1ZSVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgWmhpZ3VvJVpoYW5nJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBIYXJhbGQlUmVueiVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgWGlhbnNoZW5nJUxpdSVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgSnVuZ2FuZyVYaWUlTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIE1pbiVYaWUlTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIEppYW5waW5nJVpoYW8lTlVMTCUwXQ=="
Set-CellFromBase64 $ws 11 "W0FudXJhZGhhJUxhbGElTlVMTCUxLCAgICAgICAgICAgICAgICAgICAgICAgICAgIEtpcHAgVy4lSm9obnNvbiVOVUxMJTEsICAgICAgICAgICAgICAgICAgICAgICAgICAgSmFtZXMgTC4lSmFudXp6aSVOVUxMJTEsICAgICAgICAgICAgICAgICAgICAgICAgICAgQWRhbSBKLiVSdXNzYWslTlVMTCUxLCAgICAgICAgICAgICAgICAgICAgICAgICAgIElzaGFuJVBhcmFuanBlJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBGZWxpeCVSaWNodGVyJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBTaGFuJVpoYW8lTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIFN1bGFpbWFuJVNvbWFuaSVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgVGllbG1hbiVWYW4gVmxlY2slTlVMTCUxLCAgICAgICAgICAgICAgICAgICAgICAgICAgIEFraGlsJVZhaWQlTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIEZheXphbiVDaGF1ZGhyeSVOVUxMJTEsICAgICAgICAgICAgICAgICAgICAgICAgICAgSmVzc2ljYSBLLiVEZSBGcmVpdGFzJU5VTEwlMSwgICAgICAgICAgICAgICAgICAgICAgICAgICBaYWhpIEEuJUZheWFkJU5VTEwlMSwgICAgICAgICAgICAgICAgICAgICAgICAgICBTZWFuIFAuJVBpbm5leSVOVUxMJTEsICAgICAgICAgICAgICAgICAgICAgICAgICAgTWF0dGhldyVMZXZpbiVOVUxMJTEsICAgICAgICAgICAgICAgICAgICAgICAgICAgQWxleGFuZGVyJUNoYXJuZXklTlVMTCUxLCAgICAgICAgICAgICAgICAgICAgICAgICAgRW1pbGlhJUJhZ2llbGxhJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBKYWdhdCVOYXJ1bGElTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIEJlbmphbWluIFMuJUdsaWNrc2JlcmclTlVMTCUxLCAgICAgICAgICAgICAgICAgICAgICAgICAgIEdpcmlzaCVOYWRrYXJuaSVOVUxMJTEsICAgICAgICAgICAgICAgICAgICAgICAgICAgRG9ubmEgTS4lTWFuY2luaSVOVUxMJTEsICAgICAgICAgICAgICAgICAgICAgICAgICAgVmFsZW50aW4lRnVzdGVyJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBOVUxMJU5VTEwlTlVMTCUwXQ=="
Set-CellFromBase64 $ws 12 "W0NoYW9saW4lSHVhbmclTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIFllbWluZyVXYW5nJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBYaW5nd2FuZyVMaSVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgTGlsaSVSZW4lTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIEppYW5waW5nJVpoYW8lTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIFlpJUh1JU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBMaSVaaGFuZyVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgR3VvaHVpJUZhbiVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgSml1eWFuZyVYdSVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgWGlhb3lpbmclR3UlTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIFpoZW5zaHVuJUNoZW5nJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBUaW5nJVl1JU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBKaWFhbiVYaWElTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIFl1YW4lV2VpJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBXZW5qdWFuJVd1JU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBYdWVsZWklWGllJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBXZW4lWWluJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBIdWklTGklTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIE1pbiVMaXUlTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIFlhbiVYaWFvJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBIb25nJUdhbyVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgTGklR3VvJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBKdW5nYW5nJVhpZSVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgR3VhbmdmYSVXYW5nJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBSb25nbWVuZyVKaWFuZyVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgWmhhbmNoZW5nJUdhbyVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgUWklSmluJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBKaWFud2VpJVdhbmcld2FuZ2p3MjhAMTYzLmNvbSUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIEJpbiVDYW8lY2FvYmluX2JlbkAxNjMuY29tJTBd"
Set-CellFromBase64 $ws 13 "W0h1YW4lSGFuJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBMaW5saW4lWGllJU5VTEwlMSwgICAgICAgICAgICAgICAgICAgICAgICAgICBSdWklTGl1JU5VTEwlMSwgICAgICAgICAgICAgICAgICAgICAgICAgICBKaWUlWWFuZyVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgRmFuZyVMaXUlTlVMTCUyLCAgICAgICAgICAgICAgICAgICAgICAgICAgIEthaWxhbmclV3UlTlVMTCUyLCAgICAgICAgICAgICAgICAgICAgICAgICAgIExhbmclQ2hlbiVOVUxMJTEsICAgICAgICAgICAgICAgICAgICAgICAgICAgV2VpJUhvdSVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgWW9uZyVGZW5nJXlvbmdmZW5nQHdodS5lZHUuY24lMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBDaGVuZ2xpYW5nJVpodSV4aW5jaGVuZ3podUAxNjMuY29tJTIsICAgICAgICAgICAgICAgICAgICAgICAgICAgQ2hlbmdsaWFuZyVaaHUleGluY2hlbmd6aHVAMTYzLmNvbSUwXQ=="
Set-CellFromBase64 $ws 14 "W1FpbmclRGVuZyVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgQm8lSHUlTlVMTCUxLCAgICAgICAgICAgICAgICAgICAgICAgICAgIFlhbyVaaGFuZyVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgSGFvJVdhbmclTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIFhpYW95YW5nJVpob3UlTlVMTCUxLCAgICAgICAgICAgICAgICAgICAgICAgICAgIFdlaSVIdSVOVUxMJTEsICAgICAgICAgICAgICAgICAgICAgICAgICAgWXV0aW5nJUNoZW5nJU5VTEwlMSwgICAgICAgICAgICAgICAgICAgICAgICAgICBKaWUlWWFuJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBIYWlxaW4lUGluZyVOVUxMJTEsICAgICAgICAgICAgICAgICAgICAgICAgICAgUWluZyVaaG91JXFpbmd6aG91LndoLmVkdUBob3RtYWlsLmNvbSUxXQ=="
Set-CellFromBase64 $ws 15 "W1Jvbmdyb25nJVlhbmclTlVMTCUxLCAgICAgICAgICAgICAgICAgICAgICAgICAgIFhpZW4lR3VpJU5VTEwlMSwgICAgICAgICAgICAgICAgICAgICAgICAgICBZb25neGklWmhhbmclTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIFlvbmclWGlvbmclTlVMTCUwXQ=="
Set-CellFromBase64 $ws 16 "W1hpYW9ibyVZYW5nJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBZdWFuJVl1JU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBKaXFpYW4lWHUlTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIEh1YXFpbmclU2h1JU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBKaWEnYW4lWGlhJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBIb25nJUxpdSVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgWW9uZ3JhbiVXdSVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgTHUlWmhhbmclTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIFpodWklWXUlTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIE1pbmdoYW8lRmFuZyVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgVGluZyVZdSVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgWWF4aW4lV2FuZyVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgU2hhbmd3ZW4lUGFuJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBYaWFvamluZyVab3UlTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIFNoaXlpbmclWXVhbiVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgWW91JVNoYW5nJU5VTEwlMF0="
Set-CellFromBase64 $ws 17 "W0FndWlsYXIlTWVkaW5hIEpvc1x1MDBlOSBNLiVjb3JlR2l2ZXNOb0VtYWlsJTAsICAgICAgICAgICAgICAgICAgICAgICAgICBDbGFybyVWYWxkXHUwMGU5cyBSYW1cdTAwZjNuJWNvcmVHaXZlc05vRW1haWwlMSwgICAgICAgICAgICAgICAgICAgICAgICAgIE9icmVnXHUwMGYzbiVTYW50b3MgQW5nZWwgRy4lY29yZUdpdmVzTm9FbWFpbCUxLCAgICAgICAgICAgICAgICAgICAgICAgICAgUHJvaFx1MDBlZGFzJU1hcnRcdTAwZWRuZXogSnVhbiVjb3JlR2l2ZXNOb0VtYWlsJTEsICAgICAgICAgICAgICAgICAgICAgICAgICBSb2RyXHUwMGVkZ3VleiVCbGFuY28gU3VpbGJlcnQlY29yZUdpdmVzTm9FbWFpbCUxXQ=="
Set-CellFromBase64 $ws 18 "W1NoYW8tRmFuZyVOaWUlTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIE1pYW8lWXUlTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIFRpYW4lWGllJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBGZW4lWWFuZyVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgSG9uZy1CbyVXYW5nJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBaaGFvLUh1aSVXYW5nJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBNaW5nJUxpJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBYaW5nLUxpJUdhbyVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgQmluZy1KaWUlTHYlTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIFNoaS1KaWElV2FuZyVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgWGlhby1CbyVaaGFuZyVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgU2hhby1MaW4lSGUlTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIFpoaS1IdWElUWl1JU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBZdS1IdWElTGlhbyVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgWmktSHVhJVpob3UlTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIFhpYW5nJUNoZW5nJU5VTEwlMF0="
Set-CellFromBase64 $ws 19 "W1RhbyVHdW8lTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIFlvbmd6aGVuJUZhbiVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgTWluZyVDaGVuJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBYaWFveWFuJVd1JU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBMaW4lWmhhbmclTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIFRhbyVIZSVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgSGFpcm9uZyVXYW5nJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBKaW5nJVdhbiVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgWGluZ2h1YW4lV2FuZyVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgWmhpYmluZyVMdSVOVUxMJTBd"
Set-CellFromBase64 $ws 20 "W1lhbiVEZW5nJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBXZWklTGl1JU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBLdWklTGl1JU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBZdWFuLVl1YW4lRmFuZyVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgSmluJVNoYW5nJU5VTEwlMSwgICAgICAgICAgICAgICAgICAgICAgICAgICBMaW5nJVpob3UlTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIEtlJVdhbmclTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIEZhbiVMZW5nJU5VTEwlMSwgICAgICAgICAgICAgICAgICAgICAgICAgICBTaHVhbmclV2VpJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBMZWklQ2hlbiVOVUxMJTEsICAgICAgICAgICAgICAgICAgICAgICAgICAgSHVpLUd1byVMaXUlTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIFBlaS1GYW5nJVdlaSVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgUGVpLUZhbmclV2VpJU5VTEwlMF0="
Set-CellFromBase64 $ws 21 "W1RhbyVDaGVuJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBEaSVXdSVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgSHVpbG9uZyVDaGVuJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBXZWltaW5nJVlhbiVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgRGFubGVpJVlhbmclTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIEd1YW5nJUNoZW4lTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIEtlJU1hJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBEb25nJVh1JU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBIYWlqaW5nJVl1JU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBIb25nd3UlV2FuZyVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgVGFvJVdhbmclTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIFdlaSVHdW8lTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIEppYSVDaGVuJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBDaGVuJURpbmclTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIFhpYW9waW5nJVpoYW5nJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBKaWFxdWFuJUh1YW5nJU5VTEwlMCwgICAgICAgICAgICAgICAgICAgICAgICAgICBNZWlmYW5nJUhhbiVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgU2h1c2hlbmclTGklTlVMTCUwLCAgICAgICAgICAgICAgICAgICAgICAgICAgIFhpYW9waW5nJUx1byVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgSmlhbnBpbmclWmhhbyVOVUxMJTAsICAgICAgICAgICAgICAgICAgICAgICAgICAgUWluJU5pbmclTlVMTCUwXQ=="
